$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Felipe" row (row 5) entirely
$ws.Range("A5:F5").Delete()

# Delete the two oldest week columns (B: 11_02_2024, C: 18_02_2024)
$ws.Range("B1:C1").EntireColumn.Delete()

# Fix Betty's 11_03_2024 value (now column D after the shift)
$ws.Range("D4").Value = 1383

# Add the new week column (E) with header and values
$ws.Range("E1").Value = "17_03_2024"
$ws.Range("E2").Value = 1487
$ws.Range("E3").Value = 1448
$ws.Range("E4").Value = 1412
$ws.Range("E5").Value = 261

# Leave the same selection state the author ended up with (columns B:C selected)
$null = $ws.Range("B1:C1048576").Select()
